$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,13

$arr[0,0] = 1.04385182891744
$arr[0,1] = 0.09413454237079577
$arr[0,2] = 0.04303516009729691
$arr[0,3] = 0.03908907544184004
$arr[0,4] = 1.329879858142661
$arr[0,5] = 1.07241976931067
$arr[0,6] = 0.02000740078472996
$arr[0,7] = 0
$arr[0,8] = 0.7570710835671548
$arr[0,9] = 0.8168144458957656
$arr[0,10] = 0.05747954986053472
$arr[0,11] = 0.9702809813241231
$arr[0,12] = 0.199986060688687

$arr[1,0] = 0.9098738871720684
$arr[1,1] = 0.08344099540249772
$arr[1,2] = 0.03767649629832448
$arr[1,3] = 0.03558659694449595
$arr[1,4] = 1.283076519962435
$arr[1,5] = 1.035834600269681
$arr[1,6] = 0.02499162089995866
$arr[1,7] = 0
$arr[1,8] = 0.7457032106127315
$arr[1,9] = 0.8041743023258547
$arr[1,10] = 0.053755650805833
$arr[1,11] = 0.8434838215531499
$arr[1,12] = 0.1764159028750285

$arr[2,0] = 0.8278203139691698
$arr[2,1] = 0.07688135669327067
$arr[2,2] = 0.03438998101803747
$arr[2,3] = 0.0334355759026419
$arr[2,4] = 1.255550508324681
$arr[2,5] = 1.014525289958144
$arr[2,6] = 0.02844236893317553
$arr[2,7] = 0
$arr[2,8] = 0.7394000595923984
$arr[2,9] = 0.7971696380993123
$arr[2,10] = 0.05145364667789565
$arr[2,11] = 0.7658088755837866
$arr[2,12] = 0.1619586368864603

$arr[3,0] = 0.7944247210909055
$arr[3,1] = 0.07420829776452109
$arr[3,2] = 0.03305120043116716
$arr[3,3] = 0.03255835576293187
$arr[3,4] = 1.244624346282848
$arr[3,5] = 1.006118045731782
$arr[3,6] = 0.02994306559709115
$arr[3,7] = 0
$arr[3,8] = 0.7369954968174568
$arr[3,9] = 0.7944980436141265
$arr[3,10] = 0.05051093943649398
$arr[3,11] = 0.7341897290415034
$arr[3,12] = 0.156068832352517

$arr[4,0] = 0.7888816367531319
$arr[4,1] = 0.07376439834106918
$arr[4,2] = 0.03282891502155394
$arr[4,3] = 0.03241263878487821
$arr[4,4] = 1.242827246686744
$arr[4,5] = 1.004738337657344
$arr[4,6] = 0.03019785707338918
$arr[4,7] = 0
$arr[4,8] = 0.7366059529552729
$arr[4,9] = 0.794065263248612
$arr[4,10] = 0.05035410396093809
$arr[4,11] = 0.7289411396147045
$arr[4,12] = 0.1550908750192619

$arr[5,0] = 0.8273697734228165
$arr[5,1] = 0.07684530864536043
$arr[5,2] = 0.03437192428091151
$arr[5,3] = 0.03342374877847654
$arr[5,4] = 1.255401994058118
$arr[5,5] = 1.014410805196022
$arr[5,6] = 0.0284622300320736
$arr[5,7] = 0
$arr[5,8] = 0.7393669746237208
$arr[5,9] = 0.7971328771058523
$arr[5,10] = 0.05144095264333615
$arr[5,11] = 0.7653823250238929
$arr[5,12] = 0.1618792009460535

$arr[6,0] = 0.9976060615195763
$arr[6,1] = 0.09044516546074988
$arr[6,2] = 0.04118644643164515
$arr[6,3] = 0.03788118777471716
$arr[6,4] = 1.313482172169728
$arr[6,5] = 1.059557182566778
$arr[6,6] = 0.02164266604765264
$arr[6,7] = 0
$arr[6,8] = 0.7530071783174606
$arr[6,9] = 0.8122945955001626
$arr[6,10] = 0.05619832601198027
$arr[6,11] = 0.9265173155959587
$arr[6,12] = 0.1918546735467288

$arr[7,0] = 1.333633087740225
$arr[7,1] = 0.117238240734423
$arr[7,2] = 0.05460107383278512
$arr[7,3] = 0.046644357997355
$arr[7,4] = 1.43764552462531
$arr[7,5] = 1.157903150517811
$arr[7,6] = 0.01154045504980085
$arr[7,7] = 0
$arr[7,8] = 0.7854194782169515
$arr[7,9] = 0.8483787293206788
$arr[7,10] = 0.06543866836499035
$arr[7,11] = 1.244465600278176
$arr[7,12] = 0.2508573580104638

$arr[8,0] = 1.581026679759077
$arr[8,1] = 0.1382524548818935
$arr[8,2] = 0.06385077838980635
$arr[8,3] = 0.05124914067116215
$arr[8,4] = 1.512600742768797
$arr[8,5] = 1.216307131898077
$arr[8,6] = 0.006756977203231873
$arr[8,7] = 0
$arr[8,8] = 0.802629599055507
$arr[8,9] = 0.8653861252011836
$arr[8,10] = 0.06956319995316207
$arr[8,11] = 1.481371570855259
$arr[8,12] = 0.2857524319889251

$arr[9,0] = 1.679676918844649
$arr[9,1] = 0.158399345217191
$arr[9,2] = 0.06226170435826361
$arr[9,3] = 0.03889697728956598
$arr[9,4] = 1.350760895222507
$arr[9,5] = 1.070278644918787
$arr[9,6] = 0.02494256150944452
$arr[9,7] = 0
$arr[9,8] = 0.7233311183641149
$arr[9,9] = 0.7589335919924025
$arr[9,10] = 0.05384546461752215
$arr[9,11] = 1.601337501490207
$arr[9,12] = 0.2300614652554316

$arr[10,0] = 1.710085805698583
$arr[10,1] = 0.1717683986943968
$arr[10,2] = 0.05896221747661201
$arr[10,3] = 0.03007609817221812
$arr[10,4] = 1.205631148570987
$arr[10,5] = 0.9420391891068505
$arr[10,6] = 0.06382607657661765
$arr[10,7] = 0
$arr[10,8] = 0.6561419178488705
$arr[10,9] = 0.6718870886679156
$arr[10,10] = 0.04716832261739123
$arr[10,11] = 1.652968388415559
$arr[10,12] = 0.1814280867860489

$arr[11,0] = 1.689404929066853
$arr[11,1] = 0.180524655002344
$arr[11,2] = 0.05417864124206773
$arr[11,3] = 0.02348118466613869
$arr[11,4] = 1.064817548443614
$arr[11,5] = 0.8197901875716553
$arr[11,6] = 0.1203077799619194
$arr[11,7] = 0
$arr[11,8] = 0.5941504679873475
$arr[11,9] = 0.5944281606999766
$arr[11,10] = 0.0469475144159972
$arr[11,11] = 1.654431021794124
$arr[11,12] = 0.1357477896342019

$arr[12,0] = 1.651881221375248
$arr[12,1] = 0.1845857704195168
$arr[12,2] = 0.05019513139547627
$arr[12,3] = 0.02024935570076458
$arr[12,4] = 0.9699526509482297
$arr[12,5] = 0.7385424213113794
$arr[12,6] = 0.1704115667775739
$arr[12,7] = 0
$arr[12,8] = 0.5538466026917916
$arr[12,9] = 0.5457185100955257
$arr[12,10] = 0.05064480719824616
$arr[12,11] = 1.632763139520108
$arr[12,12] = 0.1062125138322543

$arr[13,0] = 1.631945570981372
$arr[13,1] = 0.1844594421863519
$arr[13,2] = 0.04891929610958812
$arr[13,3] = 0.01960786780551249
$arr[13,4] = 0.9469575158179282
$arr[13,5] = 0.7191799074273177
$arr[13,6] = 0.1832851461355176
$arr[13,7] = 0
$arr[13,8] = 0.5447401828889298
$arr[13,9] = 0.5351090748535121
$arr[13,10] = 0.05192530446023724
$arr[13,11] = 1.616698800744473
$arr[13,12] = 0.09898643259576545

$arr[14,0] = 1.529549076367914
$arr[14,1] = 0.1731938006278142
$arr[14,2] = 0.04604350640497756
$arr[14,3] = 0.01916845045142468
$arr[14,4] = 0.9443675443237325
$arr[14,5] = 0.7194345820847161
$arr[14,6] = 0.1725239141063497
$arr[14,7] = 0
$arr[14,8] = 0.5499094720968856
$arr[14,9] = 0.543524785066996
$arr[14,10] = 0.05014090140552341
$arr[14,11] = 1.513521128751336
$arr[14,12] = 0.09478801165559503

$arr[15,0] = 1.47114824270227
$arr[15,1] = 0.162429552206504
$arr[15,2] = 0.04593055183283212
$arr[15,3] = 0.02028359774163313
$arr[15,4] = 0.9920694818564613
$arr[15,5] = 0.7621740152973757
$arr[15,6] = 0.135901724484853
$arr[15,7] = 0
$arr[15,8] = 0.5750670005019316
$arr[15,9] = 0.5753970560991242
$arr[15,10] = 0.04617597740558033
$arr[15,11] = 1.446171835268274
$arr[15,12] = 0.106454031865745

$arr[16,0] = 1.443485568980748
$arr[16,1] = 0.1511726934994897
$arr[16,2] = 0.04819784930140258
$arr[16,3] = 0.02414488059451969
$arr[16,4] = 1.092574137302357
$arr[16,5] = 0.8503483544099026
$arr[16,6] = 0.0833505223943547
$arr[16,7] = 0
$arr[16,8] = 0.6221556222961482
$arr[16,9] = 0.6344758810235902
$arr[16,10] = 0.04348286388442579
$arr[16,11] = 1.401871759684923
$arr[16,12] = 0.1361294674065405

$arr[17,0] = 1.442312312470534
$arr[17,1] = 0.140754561336891
$arr[17,2] = 0.05228102380694821
$arr[17,3] = 0.03188109699934039
$arr[17,4] = 1.235367234166063
$arr[17,5] = 0.9755414300712886
$arr[17,6] = 0.03733144570911406
$arr[17,7] = 0
$arr[17,8] = 0.686672370794227
$arr[17,9] = 0.7166821634417104
$arr[17,10] = 0.047338448181522
$arr[17,11] = 1.379459885735457
$arr[17,12] = 0.182956452347014

$arr[18,0] = 1.515833234007602
$arr[18,1] = 0.1327607553859167
$arr[18,2] = 0.06138968585533178
$arr[18,3] = 0.04995643730191013
$arr[18,4] = 1.491408842103837
$arr[18,5] = 1.199586408788662
$arr[18,6] = 0.007853433249215946
$arr[18,7] = 0
$arr[18,8] = 0.7973844356789215
$arr[18,9] = 0.8600106756626502
$arr[18,10] = 0.06836213980171735
$arr[18,11] = 1.419052624337866
$arr[18,12] = 0.2762113891483438

$arr[19,0] = 1.709090123298353
$arr[19,1] = 0.1472167768396133
$arr[19,2] = 0.06955365349595866
$arr[19,3] = 0.05643383594834717
$arr[19,4] = 1.587842515682993
$arr[19,5] = 1.279073608950966
$arr[19,6] = 0.004377885791201019
$arr[19,7] = 0
$arr[19,8] = 0.8281827106783766
$arr[19,9] = 0.8961009913993365
$arr[19,10] = 0.07566050201043595
$arr[19,11] = 1.599718217980268
$arr[19,12] = 0.3166317263185334

$arr[20,0] = 1.835235696256291
$arr[20,1] = 0.1573160258188722
$arr[20,2] = 0.07457196335742822
$arr[20,3] = 0.05972965512117057
$arr[20,4] = 1.640542701797727
$arr[20,5] = 1.322016913624765
$arr[20,6] = 0.00279065913575649
$arr[20,7] = 0
$arr[20,8] = 0.8438180738465348
$arr[20,9] = 0.9135782301818693
$arr[20,10] = 0.07908833102521129
$arr[20,11] = 1.71909888397056
$arr[20,12] = 0.3387085952494573

$arr[21,0] = 1.767842986377332
$arr[21,1] = 0.1519184373835145
$arr[21,2] = 0.07189121680414701
$arr[21,3] = 0.05796828901906181
$arr[21,4] = 1.612258861819555
$arr[21,5] = 1.298945034679306
$arr[21,6] = 0.003590990724831822
$arr[21,7] = 0
$arr[21,8] = 0.835392102578993
$arr[21,9] = 0.9041578803150188
$arr[21,10] = 0.07725702876428997
$arr[21,11] = 1.655318256562481
$arr[21,12] = 0.3269151904721923

$arr[22,0] = 1.513963111843537
$arr[22,1] = 0.1316246593157899
$arr[22,2] = 0.06178629035027683
$arr[22,3] = 0.05134363427475641
$arr[22,4] = 1.508460497977012
$arr[22,5] = 1.214785821837438
$arr[22,6] = 0.00761841127687346
$arr[22,7] = 0
$arr[22,8] = 0.8052062991225597
$arr[22,9] = 0.8704450433922233
$arr[22,10] = 0.07035471921665248
$arr[22,11] = 1.415082707438671
$arr[22,12] = 0.2824628892708034

$arr[23,0] = 1.242413470630026
$arr[23,1] = 0.1099656713310111
$arr[23,2] = 0.05096325016649672
$arr[23,3] = 0.04426746690341687
$arr[23,4] = 1.402860046063779
$arr[23,5] = 1.13015071087662
$arr[23,6] = 0.01389805101506936
$arr[23,7] = 0
$arr[23,8] = 0.7760014960493464
$arr[23,9] = 0.8378855945724766
$arr[23,10] = 0.06294339378322533
$arr[23,11] = 1.158160173334039
$arr[23,12] = 0.2348565329258889

$ws.Range("B2:N25").Value = $arr
Write-Output "done"